$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1543.1666
$ws.Range("J43").Value = 1611.7
$ws.Range("L43").Value = 1611.7
$ws.Range("N43").Value = -1749.7
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 1264.826
$ws.Range("I132").Value = 1254.1364
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 3762.4092
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -1232.4092
$ws.Range("N132").Value = -9560
$ws.Range("H137").Value = 1628.5714
$ws.Range("I137").Value = 1345.4546
$ws.Range("J137").Value = 2666.6667
$ws.Range("K137").Value = 4036.3638
$ws.Range("L137").Value = 8000.000100000001
$ws.Range("M137").Value = -1486.3638
$ws.Range("N137").Value = -13100.0001
$ws.Range("H139").Value = 48156.715
$ws.Range("J139").Value = 48156.715
$ws.Range("L139").Value = 48156.715
$ws.Range("N139").Value = -58436.715

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3323014.2
$ws.Range("J2").Value = 998.5
$ws.Range("L2").Value = 998.5
$ws.Range("N2").Value = -1224.5
$ws.Range("H32").Value = 3906.48
$ws.Range("I32").Value = 2392.5715
$ws.Range("K32").Value = 2392.5715
$ws.Range("M32").Value = -2105.5715
$ws.Range("H74").Value = 1179.7333
$ws.Range("I74").Value = 924.6316
$ws.Range("K74").Value = 924.6316
$ws.Range("M74").Value = -50.63160000000005
$ws.Range("H77").Value = 1179.7333
$ws.Range("I77").Value = 924.6316
$ws.Range("K77").Value = 4623.158
$ws.Range("M77").Value = -255.1580000000004
$ws.Range("H110").Value = 1857.0454
$ws.Range("I110").Value = 1354.7142
$ws.Range("J110").Value = 2736.125
$ws.Range("K110").Value = 1354.7142
$ws.Range("L110").Value = 2736.125
$ws.Range("M110").Value = 690.2858000000001
$ws.Range("N110").Value = -6826.125
$ws.Range("H116").Value = 3323014.2
$ws.Range("J116").Value = 998.5
$ws.Range("L116").Value = 998.5
$ws.Range("N116").Value = -5586.5
$ws.Range("H132").Value = 1291.6863
$ws.Range("I132").Value = 922.3214
$ws.Range("K132").Value = 2766.9642
$ws.Range("M132").Value = -236.9642000000003

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3323014.2
$ws.Range("J3").Value = 998.5
$ws.Range("L3").Value = 998.5
$ws.Range("N3").Value = -1226.5
$ws.Range("H86").Value = 80175.89
$ws.Range("I86").Value = 1498.6471
$ws.Range("J86").Value = 201768
$ws.Range("K86").Value = 1498.6471
$ws.Range("L86").Value = 201768
$ws.Range("M86").Value = -375.6470999999999
$ws.Range("N86").Value = -204014
$ws.Range("H89").Value = 80175.89
$ws.Range("I89").Value = 1498.6471
$ws.Range("J89").Value = 201768
$ws.Range("K89").Value = 7493.2355
$ws.Range("L89").Value = 1008840
$ws.Range("M89").Value = -1877.2355
$ws.Range("N89").Value = -1020072
$ws.Range("H107").Value = 1934.375
$ws.Range("I107").Value = 1995.8334
$ws.Range("K107").Value = 1995.8334
$ws.Range("M107").Value = -75.83339999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1450
$ws.Range("J22").Value = 1833.3334
$ws.Range("L22").Value = 1833.3334
$ws.Range("N22").Value = -2533.3334
$ws.Range("H53").Value = 63000
$ws.Range("J53").Value = 63000
$ws.Range("L53").Value = 63000
$ws.Range("N53").Value = -64214
$ws.Range("H94").Value = 1319.8
$ws.Range("I94").Value = 1275
$ws.Range("K94").Value = 1275
$ws.Range("M94").Value = -824
$ws.Range("H134").Value = 1010.8889
$ws.Range("I134").Value = 1014
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 3042
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -507
$ws.Range("N134").Value = -8070

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11703.344
$ws.Range("J131").Value = 13066.491
$ws.Range("L131").Value = 39199.473
$ws.Range("N131").Value = -49279.473
$ws.Range("H132").Value = 924.75
$ws.Range("J132").Value = 1350
$ws.Range("L132").Value = 12150
$ws.Range("N132").Value = -17210
$ws.Range("H134").Value = 2695.0625
$ws.Range("I134").Value = 2102.9092
$ws.Range("J134").Value = 3997.8
$ws.Range("K134").Value = 6308.7276
$ws.Range("L134").Value = 11993.4
$ws.Range("M134").Value = -1238.7276
$ws.Range("N134").Value = -22133.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("H122").Value = 2056.4546
$ws.Range("I122").Value = 1737.6666
$ws.Range("K122").Value = 5212.9998
$ws.Range("M122").Value = -2762.9998
$ws.Range("H126").Value = 1769881.1
$ws.Range("I126").Value = 2528020.2
$ws.Range("J126").Value = 101974.9
$ws.Range("K126").Value = 7584060.600000001
$ws.Range("L126").Value = 305924.7
$ws.Range("M126").Value = -7581590.600000001
$ws.Range("N126").Value = -310864.7
$ws.Range("H132").Value = 1133676.6
$ws.Range("I132").Value = 1426269.1
$ws.Range("K132").Value = 4278807.300000001
$ws.Range("M132").Value = -4276277.300000001
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120
$ws.Range("H135").Value = 27999.5
$ws.Range("J135").Value = 27999.5
$ws.Range("L135").Value = 27999.5
$ws.Range("N135").Value = -38139.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6712.25
$ws.Range("I40").Value = 1750
$ws.Range("K40").Value = 1750
$ws.Range("M40").Value = -1614
$ws.Range("H122").Value = 6312
$ws.Range("J122").Value = 6875.5
$ws.Range("L122").Value = 20626.5
$ws.Range("N122").Value = -25526.5
$ws.Range("H132").Value = 1769.6177
$ws.Range("I132").Value = 1531
$ws.Range("K132").Value = 4593
$ws.Range("M132").Value = -2063
$ws.Range("H136").Value = 2772.3809
$ws.Range("I136").Value = 1338.1818
$ws.Range("K136").Value = 4014.5454
$ws.Range("M136").Value = -1464.5454

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1739
$ws.Range("I81").Value = 1983.3334
$ws.Range("K81").Value = 3966.6668
$ws.Range("M81").Value = -2905.6668
$ws.Range("H84").Value = 1739
$ws.Range("I84").Value = 1983.3334
$ws.Range("K84").Value = 19833.334
$ws.Range("M84").Value = -14529.334
$ws.Range("H107").Value = 623.9524
$ws.Range("J107").Value = 907.3333
$ws.Range("L107").Value = 2721.9999
$ws.Range("N107").Value = -6561.9999
$ws.Range("H122").Value = 157895.83
$ws.Range("I122").Value = 188975
$ws.Range("K122").Value = 566925
$ws.Range("M122").Value = -564475
$ws.Range("H126").Value = 8699.467000000001
$ws.Range("J126").Value = 9498.166999999999
$ws.Range("L126").Value = 28494.501
$ws.Range("N126").Value = -33434.501
$ws.Range("H132").Value = 2061
$ws.Range("I132").Value = 1645.0588
$ws.Range("K132").Value = 4935.1764
$ws.Range("M132").Value = -2405.1764
